$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell K1: a real date/time value (2018-02-16 12:17:33) formatted as
# "yyyy-mm-dd h:mm:ss" -- this introduces the custom number format (numFmtId
# 164) and the second cellXfs style that the sheet's dimension/spans grow to
# include.
$ws.Range("K1").Value = 43147.51219534304
$ws.Range("K1").NumberFormat = "yyyy-mm-dd h:mm:ss"

# New cell Z1: same text as A1/P1 (reuses the existing shared string rather
# than creating a new one).
$text = "`n____________<TK>____________ (@tlgkyck) • Instagram photos and videos`n"
$ws.Range("Z1").Value = $text

# Writing that multi-line text bumps row 1's height via autofit as a side
# effect; restore the row to its natural (non-custom) height so row 1 is
# left unchanged otherwise.
$ws.Rows.Item(1).EntireRow.AutoFit() | Out-Null
